$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2614.4443
$ws.Range("I32").Value = 1165.3334
$ws.Range("J32").Value = 3339
$ws.Range("K32").Value = 1165.3334
$ws.Range("L32").Value = 3339
$ws.Range("M32").Value = -839.3334
$ws.Range("N32").Value = -3991
$ws.Range("H86").Value = 86808260
$ws.Range("I86").Value = 95240110
$ws.Range("K86").Value = 95240110
$ws.Range("M86").Value = -95238987
$ws.Range("H89").Value = 86808260
$ws.Range("I89").Value = 95240110
$ws.Range("K89").Value = 476200550
$ws.Range("M89").Value = -476194934
$ws.Range("H132").Value = 1796.8182
$ws.Range("I132").Value = 1421.4642
$ws.Range("K132").Value = 4264.392599999999
$ws.Range("M132").Value = -1734.392599999999
$ws.Range("H135").Value = 556102.5
$ws.Range("I135").Value = 588732.0600000001
$ws.Range("J135").Value = 1400
$ws.Range("K135").Value = 5298588.540000001
$ws.Range("L135").Value = 12600
$ws.Range("M135").Value = -5296053.540000001
$ws.Range("N135").Value = -17670

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 47813.777
$ws.Range("I74").Value = 64577.58
$ws.Range("J74").Value = 4227.9
$ws.Range("K74").Value = 64577.58
$ws.Range("L74").Value = 4227.9
$ws.Range("M74").Value = -63703.58
$ws.Range("N74").Value = -5975.9
$ws.Range("H77").Value = 47813.777
$ws.Range("I77").Value = 64577.58
$ws.Range("J77").Value = 4227.9
$ws.Range("K77").Value = 322887.9
$ws.Range("L77").Value = 21139.5
$ws.Range("M77").Value = -318519.9
$ws.Range("N77").Value = -29875.5
$ws.Range("H97").Value = 19508.334
$ws.Range("I97").Value = 1384.75
$ws.Range("J97").Value = 55755.5
$ws.Range("K97").Value = 1384.75
$ws.Range("L97").Value = 55755.5
$ws.Range("M97").Value = -888.75
$ws.Range("N97").Value = -56747.5
$ws.Range("H132").Value = 5741.8037
$ws.Range("I132").Value = 4254.4443
$ws.Range("K132").Value = 12763.3329
$ws.Range("M132").Value = -10233.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6273.6665
$ws.Range("I134").Value = 2575.85
$ws.Range("J134").Value = 10166.105
$ws.Range("K134").Value = 7727.549999999999
$ws.Range("L134").Value = 30498.315
$ws.Range("M134").Value = -5192.549999999999
$ws.Range("N134").Value = -35568.315

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3814.5356
$ws.Range("I16").Value = 2150.5293
$ws.Range("J16").Value = 6386.1816
$ws.Range("K16").Value = 2150.5293
$ws.Range("L16").Value = 6386.1816
$ws.Range("M16").Value = -1863.5293
$ws.Range("N16").Value = -6960.1816
$ws.Range("H31").Value = 11081.275
$ws.Range("I31").Value = 3892.125
$ws.Range("J31").Value = 13820
$ws.Range("K31").Value = 3892.125
$ws.Range("L31").Value = 13820
$ws.Range("M31").Value = -3597.125
$ws.Range("N31").Value = -14410
$ws.Range("H34").Value = 11081.275
$ws.Range("I34").Value = 3892.125
$ws.Range("J34").Value = 13820
$ws.Range("K34").Value = 3892.125
$ws.Range("L34").Value = 13820
$ws.Range("M34").Value = -3690.125
$ws.Range("N34").Value = -14224
$ws.Range("H113").Value = 3814.5356
$ws.Range("I113").Value = 2150.5293
$ws.Range("J113").Value = 6386.1816
$ws.Range("K113").Value = 2150.5293
$ws.Range("L113").Value = 6386.1816
$ws.Range("M113").Value = 19.47069999999985
$ws.Range("N113").Value = -10726.1816
$ws.Range("H132").Value = 6505.533
$ws.Range("I132").Value = 3047.8333
$ws.Range("K132").Value = 9143.499899999999
$ws.Range("M132").Value = -6613.499899999999
$ws.Range("H134").Value = 8900.360000000001
$ws.Range("I134").Value = 4573.75
$ws.Range("K134").Value = 13721.25
$ws.Range("M134").Value = -11186.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 6642.6665
$ws.Range("I56").Value = 6642.6665
$ws.Range("K56").Value = 6642.6665
$ws.Range("M56").Value = -6112.6665
$ws.Range("H68").Value = 3815.75
$ws.Range("J68").Value = 5185.2
$ws.Range("L68").Value = 15555.6
$ws.Range("N68").Value = -17177.6
$ws.Range("H71").Value = 3815.75
$ws.Range("J71").Value = 5185.2
$ws.Range("L71").Value = 46666.8
$ws.Range("N71").Value = -54778.8
$ws.Range("H119").Value = 5556.75
$ws.Range("I119").Value = 5556.75
$ws.Range("K119").Value = 16670.25
$ws.Range("M119").Value = -11832.25
$ws.Range("H126").Value = 1633
$ws.Range("I126").Value = 1633
$ws.Range("K126").Value = 4899
$ws.Range("M126").Value = 41

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 17499.5
$ws.Range("I35").Value = 4999
$ws.Range("K35").Value = 4999
$ws.Range("M35").Value = -4701
$ws.Range("H70").Value = 8420.912
$ws.Range("I70").Value = 7289.5264
$ws.Range("K70").Value = 7289.5264
$ws.Range("M70").Value = -7019.5264
$ws.Range("H73").Value = 8420.912
$ws.Range("I73").Value = 7289.5264
$ws.Range("K73").Value = 7289.5264
$ws.Range("M73").Value = -6353.5264
$ws.Range("H97").Value = 6299.6
$ws.Range("I97").Value = 5749
$ws.Range("K97").Value = 5749
$ws.Range("M97").Value = -5253
$ws.Range("H102").Value = 4949.8887
$ws.Range("I102").Value = 3983.5454
$ws.Range("K102").Value = 3983.5454
$ws.Range("M102").Value = -2361.5454
$ws.Range("H113").Value = 5875.2856
$ws.Range("J113").Value = 8371.087
$ws.Range("L113").Value = 8371.087
$ws.Range("N113").Value = -12711.087
$ws.Range("H132").Value = 5274.448
$ws.Range("I132").Value = 2065.8096
$ws.Range("K132").Value = 6197.4288
$ws.Range("M132").Value = -3667.4288

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6959.7393
$ws.Range("I7").Value = 6159
$ws.Range("K7").Value = 6159
$ws.Range("M7").Value = -6047
$ws.Range("H61").Value = 2634163.8
$ws.Range("I61").Value = 3572308
$ws.Range("K61").Value = 3572308
$ws.Range("M61").Value = -3572106
$ws.Range("H92").Value = 53942
$ws.Range("J92").Value = 53942
$ws.Range("L92").Value = 53942
$ws.Range("N92").Value = -58934
$ws.Range("H113").Value = 2634163.8
$ws.Range("I113").Value = 3572308
$ws.Range("K113").Value = 3572308
$ws.Range("M113").Value = -3570138
$ws.Range("H126").Value = 6959.7393
$ws.Range("I126").Value = 6159
$ws.Range("K126").Value = 18477
$ws.Range("M126").Value = -16007
$ws.Range("H132").Value = 11911702
$ws.Range("I132").Value = 25002816
$ws.Range("J132").Value = 10688.818
$ws.Range("K132").Value = 75008448
$ws.Range("L132").Value = 32066.454
$ws.Range("M132").Value = -75005918
$ws.Range("N132").Value = -37126.454
$ws.Range("H136").Value = 8895.066000000001
$ws.Range("I136").Value = 2321
$ws.Range("K136").Value = 6963
$ws.Range("M136").Value = -4413
$ws.Range("H139").Value = 89329.336
$ws.Range("J139").Value = 89329.336
$ws.Range("L139").Value = 89329.336
$ws.Range("N139").Value = -99609.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 20850874
$ws.Range("I132").Value = 29419232
$ws.Range("J132").Value = 41999.715
$ws.Range("K132").Value = 88257696
$ws.Range("L132").Value = 125999.145
$ws.Range("M132").Value = -88255166
$ws.Range("N132").Value = -131059.145
